# Fruta / hortaliza, semanal
# Insert a new weekly record as row 11, pushing the existing rows 11-18
# down to 12-19 (Excel carries the formatting of row 11, incl. the date
# number format in column D, along with the shifted rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 11; everything below (old rows
# 11-18) shifts down to 12-19, keeping their values & formats intact.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with this week's record.
$ws.Range("A11").Value = 5
$ws.Range("B11").Value = "Macroferia Regional de Talca"
$ws.Range("C11").Value = "Maule"
$ws.Range("D11").Value = 44466
$ws.Range("D11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 300000000
$ws.Range("G11").Value = "Espárragos"
$ws.Range("H11").Value = "Verde"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 2500
$ws.Range("K11").Value = 1800
$ws.Range("L11").Value = 1800
$ws.Range("M11").Value = 1800
$ws.Range("N11").Value = "$/kilo"
$ws.Range("O11").Value = "Provincia de Linares"
$ws.Range("P11").Value = 1800
$ws.Range("Q11").Value = 1
$ws.Range("R11").Value = "Hortaliza"
